$d = $word.ActiveDocument

# Locate the "6. Conclusion" body paragraph - the one that starts with the
# old "There were some lessons learned..." sentence - without relying on a
# hard-coded paragraph index.
$hit = $d.Content
$found = $hit.Find.Execute("There were some lessons learned in this project", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the Conclusion paragraph to rewrite."
}

# Re-seat onto the owning Paragraph object (rather than the bare Find range)
# so InsertXML replaces the whole paragraph cleanly instead of splicing into
# the middle of it.
$paragraph = $hit.Paragraphs.First
$target = $paragraph.Range.Duplicate
$target.Collapse(1)  # wdCollapseStart

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Strengths: </w:t></w:r><w:r><w:t>The strengths of our design made it easier to put together, less required knowledge needed due to no database and less required software due reduced overhead.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Weaknesses: </w:t></w:r><w:r><w:t xml:space="preserve">Our design also had some weaknesses. The lack of a database while reducing overhead, also limited the scope of users and search/retrieval. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Project Suggestions: </w:t></w:r><w:r><w:t xml:space="preserve">Encompassing a database would be a big benefit to the project but would also require a re-design. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Class Suggestions: </w:t></w:r><w:r><w:t>The Class itself could have a bit more required constraints but not so much as restrict creativity of the developers. For instance, much of the rubric was unknown until the professor answered our questions in the &#8220;Ask the Professor&#8221; section.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Development History:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">The history of changes for the project can be found at our </w:t></w:r><w:r><w:t>github page here:</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>https://github.com/timwstrickland/metadata_analyzation</w:t></w:r></w:p>
'@

$target.InsertXML($newXml)

Write-Host "Conclusion section rewritten."
